# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" sheet (right after "总计", before the existing
# "2022-Q3" sheet) with the quarter's fund-holdings detail, and updates the
# "总计" (totals) roll-up sheet with a new top row summarizing the quarter,
# pushing the previous quarters' rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet right after "总计".
# ---------------------------------------------------------------------
$totalSheetTmp = $wb.Worksheets.Item(1)
$newSheetTmp = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $totalSheetTmp)
$newSheetTmp.Name = "2022-Q4"

# Re-fetch every sheet handle fresh *after* Add() -- handles captured before
# an Add() can end up pointing at the wrong sheet once the workbook's sheet
# collection has shifted.
$totalSheet = $wb.Worksheets.Item(1)
$q3Sheet    = $wb.Worksheets.Item("2022-Q3")
$q4Sheet    = $wb.Worksheets.Item("2022-Q4")

# Pull header row (B1:H1) styling + text from the existing "2022-Q3" sheet,
# which uses the same column layout.
$q3Sheet.Range("B1:H1").Copy($q4Sheet.Range("B1:H1"))

# Stamp every data row (2-7) with the same per-column styling used by the
# other quarter sheets (column A bold/centered index style, everything else
# plain) by copying the formatting of an existing data row down across all
# six rows we are about to fill in.
$q3Sheet.Range("A2:H2").Copy($q4Sheet.Range("A2:H2"))
$q3Sheet.Range("A2:H2").Copy($q4Sheet.Range("A3:H3"))
$q3Sheet.Range("A2:H2").Copy($q4Sheet.Range("A4:H4"))
$q3Sheet.Range("A2:H2").Copy($q4Sheet.Range("A5:H5"))
$q3Sheet.Range("A2:H2").Copy($q4Sheet.Range("A6:H6"))
$q3Sheet.Range("A2:H2").Copy($q4Sheet.Range("A7:H7"))

# Columns B, D, E, F, G hold numeric-looking values (fund codes with
# leading zeros, percentages, etc.) that must stay literal text -- force
# "Text" number format before writing so codes like "003166" don't get
# reinterpreted as the number 3166. (G7's "0.00" is the one exception:
# the source data stores it as a bare numeric 0, so it is left out of this
# range and written as a real number below.)
$q4Sheet.Range("B2:B7").NumberFormat = "@"
$q4Sheet.Range("D2:F7").NumberFormat = "@"
$q4Sheet.Range("G2:G6").NumberFormat = "@"

$q4Data = @(
    @("920002", "中金精选股票A",         "3.21", "93.08", "4.07", "0.1306", 1),
    @("003166", "鹏华弘嘉灵活配置混合C", "0.92", "91.65", "4.01", "0.0369", 9),
    @("003165", "鹏华弘嘉灵活配置混合A", "0.72", "91.65", "4.01", "0.0289", 9),
    @("920922", "中金精选股票C",         "0.12", "93.08", "4.07", "0.0049", 1),
    @("000892", "九泰天宝灵活配置混合A", "0.06", "94.55", "6.80", "0.0041", 1),
    @("002028", "九泰天宝灵活配置混合C", "0.00", "94.55", "6.80", 0,        1)
)

for ($i = 0; $i -lt $q4Data.Count; $i++) {
    $row = $i + 2
    $rec = $q4Data[$i]

    $q4Sheet.Cells.Item($row, 1).Value = $i
    $q4Sheet.Cells.Item($row, 2).Value = $rec[0]
    $q4Sheet.Cells.Item($row, 3).Value = $rec[1]
    $q4Sheet.Cells.Item($row, 4).Value = $rec[2]
    $q4Sheet.Cells.Item($row, 5).Value = $rec[3]
    $q4Sheet.Cells.Item($row, 6).Value = $rec[4]
    $q4Sheet.Cells.Item($row, 7).Value = $rec[5]
    $q4Sheet.Cells.Item($row, 8).Value = $rec[6]
}

# ---------------------------------------------------------------------
# 2) Update the "总计" roll-up: shift the three existing quarter rows down
#    one row (copying each row's formatting along with it) and write the
#    new 2022-Q4 summary row on top.
# ---------------------------------------------------------------------
$totalSheet.Range("A4:D4").Copy($totalSheet.Range("A5:D5"))

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2022-Q1"
$totalSheet.Range("C5").Value = 2
$totalSheet.Range("D5").Value = 0.04

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2022-Q2"
$totalSheet.Range("C4").Value = 2
$totalSheet.Range("D4").Value = 0.01

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q3"
$totalSheet.Range("C3").Value = 3
$totalSheet.Range("D3").Value = 0.21

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 6
$totalSheet.Range("D2").Value = 0.21
